$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.908.60"
$ws.Range("E2").Value = "  +3.19%  "

$ws.Range("D3").Value = "2.257.57"
$ws.Range("E3").Value = "  +2.40%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "251.86"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").Value = "0.638"
$ws.Range("E6").Value = "  +1.25%  "

$ws.Range("D7").Value = "70.98"
$ws.Range("E7").Value = "  +2.44%  "

$ws.Range("D8").Value = "0.672"
$ws.Range("E8").Value = "  +15.41%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").Value = "39.28"
$ws.Range("E10").Value = "  +3.92%  "

$ws.Range("D11").Value = "0.0973"
$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("D12").Value = "59.53"
$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("D13").Value = "7.58"
$ws.Range("E13").Value = "  +5.60%  "

$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").Value = "2.599.49"
$ws.Range("E15").Value = "  +2.63%  "

$ws.Range("D16").Value = "0.886"
$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("D17").Value = "14.81"
$ws.Range("E17").Value = "  +1.30%  "

$ws.Range("D18").Value = "2.262.58"
$ws.Range("E18").Value = "  +3.26%  "

$ws.Range("D19").Value = "42.850.35"
$ws.Range("E19").Value = "  +3.18%  "

$ws.Range("D20").Value = "0.0₃0981"
$ws.Range("E20").Value = "  +1.62%  "

$ws.Range("D21").Value = "6.28"
$ws.Range("E21").Value = "  +0.59%  "

$ws.Range("D22").Value = "73.08"
$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("D23").Value = "237.13"
$ws.Range("E23").Value = "  +1.03%  "

$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  +3.49%  "

$ws.Range("E25").Value = "  +1.49%  "

$ws.Range("D26").Value = "11.69"
$ws.Range("E26").Value = "  -2.16%  "

$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("E28").Value = "  -2.70%  "

$ws.Range("D29").Value = "3.67"
$ws.Range("E29").Value = "  -1.45%  "

$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +2.04%  "

$ws.Range("D31").Value = "167.89"

$ws.Range("D32").Value = "21.19"
$ws.Range("E32").Value = "  +1.90%  "

$ws.Range("D33").Value = "6.24"
$ws.Range("E33").Value = "  +13.06%  "

$ws.Range("D34").Value = "0.127"
$ws.Range("E34").Value = "  +6.93%  "

$ws.Range("D35").Value = "0.0770"
$ws.Range("E35").Value = "  +1.27%  "

$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("D37").Value = "28.90"
$ws.Range("E37").Value = "  +12.70%  "

$ws.Range("D38").Value = "4.71"
$ws.Range("E38").Value = "  +1.20%  "

$ws.Range("D39").Value = "4.12"
$ws.Range("E39").Value = "  -1.90%  "

$ws.Range("D40").Value = "0.0321"
$ws.Range("E40").Value = "  +6.56%  "

$ws.Range("D41").Value = "2.30"
$ws.Range("E41").Value = "  +2.95%  "

$ws.Range("D42").Value = "5.84"
$ws.Range("E42").Value = "  +1.74%  "

$ws.Range("D43").Value = "12.10"
$ws.Range("E43").Value = "  -1.25%  "

$ws.Range("D44").Value = "64.17"
$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("E45").Value = "  +2.21%  "

$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").Value = "8.92"
$ws.Range("E47").Value = "  +1.39%  "

$ws.Range("E48").Value = "  +1.07%  "

$ws.Range("E49").Value = "  -5.18%  "

$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("E51").Value = "  +0.79%  "
